# Auto-generated script to update cryptos worksheet cell values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '37.250.13'
$c.ClearFormats()
$ws.Range('E2').Value = '  +0.88%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.062.27'
$c.ClearFormats()
$ws.Range('E3').Value = '  +1.03%  '
$ws.Range('E4').Value = '  -0.37%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '249.13'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('E6').Value = '  +1.24%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '59.30'
$c.ClearFormats()
$ws.Range('E7').Value = '  +7.26%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +2.67%  '
$ws.Range('E10').Value = '  +1.74%  '
$ws.Range('E11').Value = '  +2.06%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '15.96'
$c.ClearFormats()
$ws.Range('E12').Value = '  +2.00%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.919'
$c.ClearFormats()
$ws.Range('E13').Value = '  +17.08%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '2.364.21'
$c.ClearFormats()
$ws.Range('E14').Value = '  +1.04%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '5.85'
$c.ClearFormats()
$ws.Range('E15').Value = '  +5.29%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '2.073.45'
$c.ClearFormats()
$ws.Range('E16').Value = '  +1.42%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '18.85'
$c.ClearFormats()
$ws.Range('E17').Value = '  +14.74%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '37.242.07'
$c.ClearFormats()
$ws.Range('E18').Value = '  +0.99%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '75.57'
$c.ClearFormats()
$ws.Range('E19').Value = '  +2.97%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.0₃0911'
$c.ClearFormats()
$ws.Range('E20').Value = '  +2.57%  '
$ws.Range('E21').Value = '  +4.63%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '239.51'
$c.ClearFormats()
$ws.Range('E22').Value = '  +1.99%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range('E23').Value = '  -0.08%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.48'
$c.ClearFormats()
$ws.Range('E24').Value = '  +6.26%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.23'
$c.ClearFormats()
$ws.Range('E25').Value = '  +3.39%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '9.67'
$c.ClearFormats()
$ws.Range('E26').Value = '  +7.08%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '171.77'
$c.ClearFormats()
$ws.Range('E27').Value = '  +2.86%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '20.33'
$c.ClearFormats()
$ws.Range('E28').Value = '  +3.59%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '5.54'
$c.ClearFormats()
$ws.Range('E29').Value = '  +19.66%  '
$ws.Range('E30').Value = '  +1.45%  '
$ws.Range('E31').Value = '  +6.41%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '4.82'
$c.ClearFormats()
$ws.Range('E32').Value = '  +10.82%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.0632'
$c.ClearFormats()
$ws.Range('E33').Value = '  +4.31%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '2.35'
$c.ClearFormats()
$ws.Range('E34').Value = '  +7.21%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.0885'
$c.ClearFormats()
$ws.Range('E35').Value = '  +2.06%  '
$ws.Range('E36').Value = '  -0.15%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.84'
$c.ClearFormats()
$ws.Range('E37').Value = '  +5.27%  '
$ws.Range('E38').Value = '  +1.46%  '
$ws.Range('B39').Value = 'HuobiToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '3.13'
$c.ClearFormats()
$ws.Range('E39').Value = '  -3.33%  '
$ws.Range('B40').Value = 'THORChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '5.19'
$c.ClearFormats()
$ws.Range('E40').Value = '  +7.03%  '
$ws.Range('E41').Value = '  -3.61%  '
$ws.Range('E42').Value = '  +3.87%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '101.14'
$c.ClearFormats()
$ws.Range('E43').Value = '  +6.70%  '
$ws.Range('E44').Value = '  +6.11%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '17.53'
$c.ClearFormats()
$ws.Range('E45').Value = '  +2.70%  '
$ws.Range('E46').Value = '  +1.93%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.313.84'
$c.ClearFormats()
$ws.Range('E47').Value = '  +3.49%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '3.89'
$c.ClearFormats()
$ws.Range('E48').Value = '  +22.23%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '6.99'
$c.ClearFormats()
$ws.Range('E49').Value = '  +5.48%  '
$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.88'
$c.ClearFormats()
$ws.Range('E50').Value = '  +1.51%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '2.252.25'
$c.ClearFormats()
$ws.Range('E51').Value = '  +1.36%  '
